# Regenerate save_data column G ("K") values: recompute K (was "Strike#")
# for each data row (rows 2-54 of Sheet1), writing the newly calculated
# integer values in place of the old ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 2
    7  = 2
    8  = 2
    9  = 1
    10 = 4
    11 = 0
    12 = 2
    13 = 2
    14 = 1
    15 = 2
    16 = 1
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 3
    24 = 0
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 0
    30 = 0
    31 = 3
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 2
    38 = 0
    39 = 0
    40 = 3
    41 = 1
    42 = 0
    43 = 1
    44 = 0
    45 = 4
    46 = 1
    47 = 0
    48 = 3
    49 = 4
    50 = 0
    51 = 1
    52 = 0
    53 = 0
    54 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
